$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (B, C, D, E) keep their values as plain text,
# matching the original inlineStr cell type (avoids numeric auto-conversion,
# e.g. "0.9050" or "1.085" becoming a Double).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.193.29'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.903.49'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '307.52'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5247'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3819'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +1.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07308'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.66'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +2.38%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.9050'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08106'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -3.82%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '95.63'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.359'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +1.64%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '1.824.65'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -4.68%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008663'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '14.73'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '27.234.87'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.114'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +1.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.81'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +2.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.467'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.45%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.333'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +2.35%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '149.42'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.72%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '18.25'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.743'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '116.36'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.836'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.47%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.893'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.09243'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.56%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.05078'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.7972'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.229'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.981'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.93%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.374'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.676'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +1.48%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5736'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01998'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.085'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '9.025'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.597'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.65%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '116.59'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.80%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1517'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.4894'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.003'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.18'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '38.60'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '64.09'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05960'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.37%  '
